$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Status column (G) for the rows relating to teacher/class
# admin functionality from "In Progress" to "Testing".
$rows = @(55, 56, 57, 58, 60, 61, 62)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Testing"
}

# Update the selected cell in the sheet view to G62.
$ws.Range("G62").Select()
